$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column B; existing B..E shift to C..F,
# carrying their values/formatting along automatically.
$ws.Columns("B:B").Insert()

# New column B is this week's ("Jun_26") ratings.
$ws.Range("B1").Value = "Jun_26"
for ($r = 2; $r -le 27; $r++) {
    $ws.Cells.Item($r, 2).Value = "UN"
}

# Two new analysts added to the bottom of the table, rated this week only.
$ws.Range("A28").Value = "Benchmark"
$ws.Range("B28").Value = "UN"
$ws.Range("A29").Value = "Evercore ISI"
$ws.Range("B29").Value = "UN"
